# Fill in missing "Miền" (region) values in column E with "Chưa cập nhật"
# for all data rows that currently have no value in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$value = "Chưa cập nhật"

# Contiguous blocks of rows (1-based, matching worksheet row numbers)
# that are missing the column E ("Miền") value and need it populated.
$ranges = @(
    @(23, 132),
    @(159, 162),
    @(169, 180)
)

foreach ($range in $ranges) {
    $startRow = $range[0]
    $endRow = $range[1]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Cells.Item($r, 5).Value = $value
    }
}
